$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D cells whose new value would otherwise be auto-parsed as a number by Excel;
# force them to remain plain text (matching the source inlineStr cells) by pre-setting
# the cell number format to Text ("@") before assigning the string value.
$textPriceCells = @(
    'D5',
    'D6',
    'D13',
    'D14',
    'D17',
    'D19',
    'D22',
    'D23',
    'D26',
    'D27',
    'D28',
    'D29',
    'D31',
    'D32',
    'D37',
    'D45',
    'D47',
    'D49',
    'D51'
)
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply every updated cell value (ticker prices / % changes / coin name-link-price-volume rows).
$cellValues = [ordered]@{
    'D2' = '42.833.38'
    'E2' = '  -0.53%  '
    'D3' = '2.292.85'
    'E3' = '  -0.85%  '
    'E4' = '  -0.03%  '
    'D5' = '299.65'
    'E5' = '  -0.74%  '
    'D6' = '96.70'
    'E6' = '  -2.06%  '
    'E7' = '  +0.59%  '
    'E8' = '  -0.01%  '
    'E9' = '  -3.24%  '
    'E10' = '  -0.27%  '
    'E11' = '  -0.35%  '
    'E12' = '  +0.67%  '
    'D13' = '17.70'
    'E13' = '  -0.62%  '
    'D14' = '6.75'
    'E14' = '  -1.94%  '
    'D15' = '2.647.40'
    'E15' = '  -0.95%  '
    'D16' = '2.289.65'
    'E16' = '  +0.54%  '
    'D17' = '0.774'
    'E17' = '  -1.84%  '
    'D18' = '42.752.39'
    'E18' = '  -0.52%  '
    'D19' = '12.72'
    'E19' = '  -4.15%  '
    'D20' = '0.0₃0903'
    'E20' = '  -0.50%  '
    'E21' = '  -2.36%  '
    'D22' = '67.79'
    'E22' = '  -0.43%  '
    'D23' = '240.85'
    'E23' = '  -0.14%  '
    'E24' = '  -1.52%  '
    'E25' = '  +0.06%  '
    'D26' = '2.42'
    'E26' = '  -1.40%  '
    'D27' = '4.01'
    'E27' = '  -0.45%  '
    'D28' = '25.11'
    'E28' = '  -0.01%  '
    'D29' = '165.94'
    'E29' = '  -1.86%  '
    'E30' = '  -0.93%  '
    'D31' = '9.03'
    'E31' = '  -1.53%  '
    'D32' = '32.83'
    'E32' = '  -1.81%  '
    'E33' = '  +0.08%  '
    'E34' = '  -2.68%  '
    'E35' = '  -3.44%  '
    'E36' = '  -7.07%  '
    'D37' = '2.38'
    'E37' = '  -1.12%  '
    'E38' = '  -1.42%  '
    'E39' = '  -1.70%  '
    'E40' = '  -3.55%  '
    'E41' = '  +0.21%  '
    'E42' = '  -0.97%  '
    'D43' = '2.008.70'
    'E43' = '  +0.68%  '
    'E44' = '  -2.76%  '
    'D45' = '10.09'
    'E45' = '  -0.31%  '
    'E46' = '  +1.32%  '
    'D47' = '17.08'
    'E47' = '  -2.94%  '
    'E48' = '  -2.24%  '
    'B49' = 'HuobiToken'
    'C49' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'D49' = '2.94'
    'E49' = '  -0.88%  '
    'B50' = 'RocketPoolETH'
    'C50' = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
    'D50' = '2.514.35'
    'E50' = '  -1.01%  '
    'B51' = 'MultiversX'
    'C51' = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
    'D51' = '53.05'
    'E51' = '  -3.13%  '
}
foreach ($cellRef in $cellValues.Keys) {
    $ws.Range($cellRef).Value = $cellValues[$cellRef]
}
